$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.402708292007446
$ws.Range("B1").Value = 1.784400582313538
$ws.Range("C1").Value = 1.985599756240845
$ws.Range("D1").Value = 2.279743194580078
$ws.Range("E1").Value = 2.77721095085144
